$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.673.79"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "3.274.36"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.09"
$ws.Range("E5").Value = "  +1.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.35"
$ws.Range("E6").Value = "  -1.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").Value = "  +7.40%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -3.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.72"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.402"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "3.841.87"
$ws.Range("E12").Value = "  -1.22%  "

$ws.Range("E13").Value = "  -4.31%  "

$ws.Range("D14").Value = "65.785.11"
$ws.Range("E14").Value = "  -1.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.94"
$ws.Range("E15").Value = "  -4.48%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.359.90"
$ws.Range("E16").Value = "  +1.56%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").Value = "  -2.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "426.60"
$ws.Range("E18").Value = "  -2.03%  "

$ws.Range("E19").Value = "  -4.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.36"
$ws.Range("E21").Value = "  -3.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.97"
$ws.Range("E22").Value = "  -2.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.69"
$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("D25").Value = "3.430.11"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.508"
$ws.Range("E26").Value = "  -1.38%  "

$ws.Range("E28").Value = "  -4.92%  "

$ws.Range("E29").Value = "  -2.38%  "

$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.20"
$ws.Range("E32").Value = "  -2.61%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  -3.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.57"
$ws.Range("E35").Value = "  -3.03%  "

$ws.Range("E36").Value = "  -3.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.50"
$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("E38").Value = "  -6.09%  "

$ws.Range("E39").Value = "  -3.40%  "

$ws.Range("E40").Value = "  -3.65%  "

$ws.Range("D41").Value = "2.783.82"
$ws.Range("E41").Value = "  -1.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.765"
$ws.Range("E42").Value = "  -3.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.31"
$ws.Range("E43").Value = "  -3.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.03"

$ws.Range("E45").Value = "  -2.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.92"
$ws.Range("E46").Value = "  -5.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.28"
$ws.Range("E47").Value = "  -2.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "314.03"
$ws.Range("E48").Value = "  -1.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.06"
$ws.Range("E49").Value = "  -5.34%  "

$ws.Range("E50").Value = "  -2.25%  "

$ws.Range("E51").Value = "  +5.00%  "
